{"js": "// Update the multiplication problems in the practice-sheet table.\n// Each mapping entry is the exact original cell text -> the new cell text.\nconst replacements = {\n  \"81\u00d785=\": \"57\u00d793=\",\n  \"85\u00d795=\": \"82\u00d725=\",\n  \"90\u00d789=\": \"53\u00d792=\",\n  \"35\u00d799=\": \"75\u00d786=\",\n  \"40\u00d779=\": \"11\u00d722=\",\n  \"96\u00d788=\": \"46\u00d785=\",\n  \"51\u00d773=\": \"89\u00d745=\",\n  \"77\u00d765=\": \"43\u00d718=\",\n  \"64\u00d745=\": \"63\u00d763=\",\n  \"67\u00d753=\": \"43\u00d724=\",\n  \"12\u00d783=\": \"64\u00d753=\",\n  \"36\u00d780=\": \"18\u00d770=\",\n  \"42\u00d797=\": \"63\u00d759=\",\n  \"32\u00d742=\": \"63\u00d762=\",\n  \"58\u00d729=\": \"69\u00d776=\",\n  \"61\u00d715=\": \"53\u00d745=\",\n  \"35\u00d774=\": \"86\u00d753=\",\n  \"92\u00d757=\": \"24\u00d760=\",\n  \"76\u00d720=\": \"97\u00d761=\",\n  \"34\u00d781=\": \"93\u00d730=\",\n  \"87\u00d764=\": \"65\u00d798=\",\n  \"51\u00d758=\": \"45\u00d751=\",\n  \"48\u00d796=\": \"38\u00d784=\",\n  \"35\u00d723=\": \"37\u00d772=\",\n  \"59\u00d721=\": \"12\u00d784=\",\n};\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const current = p.text;\n  if (Object.prototype.hasOwnProperty.call(replacements, current)) {\n    p.insertText(replacements[current], \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the multiplication problems in the practice-sheet table.\n# Each mapping entry is the exact original cell text -> the new cell text.\n$d = $word.ActiveDocument\n\n$replacements = [ordered]@{\n  \"81\u00d785=\" = \"57\u00d793=\"\n  \"85\u00d795=\" = \"82\u00d725=\"\n  \"90\u00d789=\" = \"53\u00d792=\"\n  \"35\u00d799=\" = \"75\u00d786=\"\n  \"40\u00d779=\" = \"11\u00d722=\"\n  \"96\u00d788=\" = \"46\u00d785=\"\n  \"51\u00d773=\" = \"89\u00d745=\"\n  \"77\u00d765=\" = \"43\u00d718=\"\n  \"64\u00d745=\" = \"63\u00d763=\"\n  \"67\u00d753=\" = \"43\u00d724=\"\n  \"12\u00d783=\" = \"64\u00d753=\"\n  \"36\u00d780=\" = \"18\u00d770=\"\n  \"42\u00d797=\" = \"63\u00d759=\"\n  \"32\u00d742=\" = \"63\u00d762=\"\n  \"58\u00d729=\" = \"69\u00d776=\"\n  \"61\u00d715=\" = \"53\u00d745=\"\n  \"35\u00d774=\" = \"86\u00d753=\"\n  \"92\u00d757=\" = \"24\u00d760=\"\n  \"76\u00d720=\" = \"97\u00d761=\"\n  \"34\u00d781=\" = \"93\u00d730=\"\n  \"87\u00d764=\" = \"65\u00d798=\"\n  \"51\u00d758=\" = \"45\u00d751=\"\n  \"48\u00d796=\" = \"38\u00d784=\"\n  \"35\u00d723=\" = \"37\u00d772=\"\n  \"59\u00d721=\" = \"12\u00d784=\"\n}\n\nforeach ($key in $replacements.Keys) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $key\n  $find.Replacement.Text = $replacements[$key]\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute(\n    $key,            # FindText\n    $true,           # MatchCase\n    $false,          # MatchWholeWord\n    $false,          # MatchWildcards\n    $false,          # MatchSoundsLike\n    $false,          # MatchAllWordForms\n    $true,           # Forward\n    1,               # Wrap (wdFindContinue)\n    $false,          # Format\n    $replacements[$key], # ReplaceWith\n    2                # Replace (wdReplaceAll)\n  ) | Out-Null\n}\n"}
